$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6416.4165
$ws.Range("I62").Value = 3666.3333
$ws.Range("K62").Value = 3666.3333
$ws.Range("M62").Value = -3042.3333
$ws.Range("H65").Value = 6416.4165
$ws.Range("I65").Value = 3666.3333
$ws.Range("K65").Value = 18331.6665
$ws.Range("M65").Value = -15211.6665
$ws.Range("H74").Value = 5800.625
$ws.Range("I74").Value = 5463.2
$ws.Range("K74").Value = 5463.2
$ws.Range("M74").Value = -4527.2
$ws.Range("H76").Value = 6031.1113
$ws.Range("I76").Value = 4446.25
$ws.Range("J76").Value = 7299
$ws.Range("K76").Value = 4446.25
$ws.Range("L76").Value = 7299
$ws.Range("M76").Value = -4131.25
$ws.Range("N76").Value = -7929
$ws.Range("H77").Value = 5800.625
$ws.Range("I77").Value = 5463.2
$ws.Range("K77").Value = 27316
$ws.Range("M77").Value = -22636
$ws.Range("H79").Value = 6031.1113
$ws.Range("I79").Value = 4446.25
$ws.Range("J79").Value = 7299
$ws.Range("K79").Value = 4446.25
$ws.Range("L79").Value = 7299
$ws.Range("M79").Value = -3354.25
$ws.Range("N79").Value = -9483
$ws.Range("H86").Value = 4551.6
$ws.Range("J86").Value = 5364.4443
$ws.Range("L86").Value = 5364.4443
$ws.Range("N86").Value = -7610.4443
$ws.Range("H88").Value = 1625.25
$ws.Range("J88").Value = 2014.25
$ws.Range("L88").Value = 2014.25
$ws.Range("N88").Value = -2826.25
$ws.Range("H89").Value = 4551.6
$ws.Range("J89").Value = 5364.4443
$ws.Range("L89").Value = 26822.2215
$ws.Range("N89").Value = -38054.2215
$ws.Range("H91").Value = 1625.25
$ws.Range("J91").Value = 2014.25
$ws.Range("L91").Value = 2014.25
$ws.Range("N91").Value = -4822.25
$ws.Range("H98").Value = 831.4545000000001
$ws.Range("J98").Value = 2495
$ws.Range("L98").Value = 2495
$ws.Range("N98").Value = -5491
$ws.Range("H122").Value = 831.4545000000001
$ws.Range("J122").Value = 2495
$ws.Range("L122").Value = 7485
$ws.Range("N122").Value = -12385
$ws.Range("H132").Value = 1306.3478
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H141").Value = 1652.125
$ws.Range("I141").Value = 962.26666
$ws.Range("K141").Value = 2886.79998
$ws.Range("M141").Value = 2293.20002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1389.6666
$ws.Range("I61").Value = 1372.0769
$ws.Range("J61").Value = 1504
$ws.Range("K61").Value = 1372.0769
$ws.Range("L61").Value = 1504
$ws.Range("M61").Value = -1160.0769
$ws.Range("N61").Value = -1928
$ws.Range("H63").Value = 6871.7856
$ws.Range("I63").Value = 6539.8
$ws.Range("J63").Value = 7056.222
$ws.Range("K63").Value = 6539.8
$ws.Range("L63").Value = 7056.222
$ws.Range("M63").Value = -5853.8
$ws.Range("N63").Value = -8428.222
$ws.Range("H66").Value = 6871.7856
$ws.Range("I66").Value = 6539.8
$ws.Range("J66").Value = 7056.222
$ws.Range("K66").Value = 32699
$ws.Range("L66").Value = 35281.11
$ws.Range("M66").Value = -29267
$ws.Range("N66").Value = -42145.11
$ws.Range("H88").Value = 251.3
$ws.Range("I88").Value = 259.22223
$ws.Range("J88").Value = 180
$ws.Range("K88").Value = 259.22223
$ws.Range("L88").Value = 180
$ws.Range("M88").Value = 146.77777
$ws.Range("N88").Value = -992
$ws.Range("H91").Value = 251.3
$ws.Range("I91").Value = 259.22223
$ws.Range("J91").Value = 180
$ws.Range("K91").Value = 259.22223
$ws.Range("L91").Value = 180
$ws.Range("M91").Value = 1144.77777
$ws.Range("N91").Value = -2988
$ws.Range("H132").Value = 1762.7435
$ws.Range("I132").Value = 1676.3334
$ws.Range("J132").Value = 2799.6667
$ws.Range("K132").Value = 5029.0002
$ws.Range("L132").Value = 8399.000100000001
$ws.Range("M132").Value = -2499.0002
$ws.Range("N132").Value = -13459.0001
$ws.Range("H136").Value = 1389.6666
$ws.Range("I136").Value = 1372.0769
$ws.Range("J136").Value = 1504
$ws.Range("K136").Value = 4116.2307
$ws.Range("L136").Value = 4512
$ws.Range("M136").Value = -1566.2307
$ws.Range("N136").Value = -9612
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 61.57143
$ws.Range("I7").Value = 48.666668
$ws.Range("J7").Value = 139
$ws.Range("K7").Value = 48.666668
$ws.Range("L7").Value = 139
$ws.Range("M7").Value = 64.333332
$ws.Range("N7").Value = -365
$ws.Range("H31").Value = 5655.25
$ws.Range("I31").Value = 1998
$ws.Range("J31").Value = 7483.875
$ws.Range("K31").Value = 1998
$ws.Range("L31").Value = 7483.875
$ws.Range("M31").Value = -1703
$ws.Range("N31").Value = -8073.875
$ws.Range("H34").Value = 5655.25
$ws.Range("I34").Value = 1998
$ws.Range("J34").Value = 7483.875
$ws.Range("K34").Value = 1998
$ws.Range("L34").Value = 7483.875
$ws.Range("M34").Value = -1796
$ws.Range("N34").Value = -7887.875
$ws.Range("H99").Value = 10609.577
$ws.Range("I99").Value = 6915.077
$ws.Range("K99").Value = 6915.077
$ws.Range("M99").Value = -5417.077
$ws.Range("H126").Value = 10609.577
$ws.Range("I126").Value = 6915.077
$ws.Range("K126").Value = 20745.231
$ws.Range("M126").Value = -18275.231
$ws.Range("H132").Value = 2567.8696
$ws.Range("I132").Value = 1716.8667
$ws.Range("K132").Value = 5150.6001
$ws.Range("M132").Value = -2620.6001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1008
$ws.Range("I14").Value = 1008
$ws.Range("K14").Value = 3024
$ws.Range("M14").Value = -2851
$ws.Range("H68").Value = 1700
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1700
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5100
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6722
$ws.Range("H71").Value = 1700
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1700
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 15300
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -23412
$ws.Range("H86").Value = 156.85715
$ws.Range("I86").Value = 102.5
$ws.Range("J86").Value = 229.33333
$ws.Range("K86").Value = 307.5
$ws.Range("L86").Value = 687.99999
$ws.Range("M86").Value = 878.5
$ws.Range("N86").Value = -3059.99999
$ws.Range("H89").Value = 156.85715
$ws.Range("I89").Value = 102.5
$ws.Range("J89").Value = 229.33333
$ws.Range("K89").Value = 922.5
$ws.Range("L89").Value = 2063.99997
$ws.Range("M89").Value = 5005.5
$ws.Range("N89").Value = -13919.99997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1846.1578
$ws.Range("I132").Value = 1454.9445
$ws.Range("K132").Value = 4364.833500000001
$ws.Range("M132").Value = -1834.833500000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1999
$ws.Range("I7").Value = 1999
$ws.Range("K7").Value = 1999
$ws.Range("M7").Value = -1887
$ws.Range("H94").Value = 35000
$ws.Range("J94").Value = 35000
$ws.Range("L94").Value = 35000
$ws.Range("N94").Value = -36352
$ws.Range("H122").Value = 10151.23
$ws.Range("I122").Value = 9747.166999999999
$ws.Range("K122").Value = 29241.501
$ws.Range("M122").Value = -26791.501
$ws.Range("H126").Value = 1999
$ws.Range("I126").Value = 1999
$ws.Range("K126").Value = 5997
$ws.Range("M126").Value = -3527
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5312.091
$ws.Range("I122").Value = 5733.4
$ws.Range("J122").Value = 1099
$ws.Range("K122").Value = 17200.2
$ws.Range("L122").Value = 3297
$ws.Range("M122").Value = -14750.2
$ws.Range("N122").Value = -8197
$ws.Range("H126").Value = 1749.5
$ws.Range("I126").Value = 1749.5
$ws.Range("K126").Value = 5248.5
$ws.Range("M126").Value = -2778.5
$ws.Range("H136").Value = 2028.9722
$ws.Range("I136").Value = 635.3077
$ws.Range("J136").Value = 5652.5
$ws.Range("K136").Value = 1905.9231
$ws.Range("L136").Value = 16957.5
$ws.Range("M136").Value = 644.0769
$ws.Range("N136").Value = -22057.5
